$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: AD1 "Wins", AE1 "Losses", AF1 "Ties" -- styled like the
# existing header row (bold, bordered, centered/top-aligned).
$hdr = $ws.Range("AD1:AF1")
$hdr.Font.Bold = $true
$hdr.HorizontalAlignment = -4108  # xlCenter
$hdr.VerticalAlignment = -4160    # xlTop
$hdr.Borders.LineStyle = 1        # xlContinuous
$hdr.Borders.Weight = 2           # xlThin

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Data rows 2-50: team record Wins/Losses/Ties for every player row.
$ws.Range("AD2:AD50").Value = 70
$ws.Range("AE2:AE50").Value = 92
$ws.Range("AF2:AF50").Value = 0

Write-Output "done"
